$wb = $excel.ActiveWorkbook

# --- Sheet "order": tighten the C-column validation messages ---------------
$orderWs = $wb.Worksheets.Item("order")

$orderWs.Range("C1").Formula  = '=IF(LEN(B1) =0,"必须选择","")'
$orderWs.Range("C2").Formula  = '=IF(LEN(B2) =0,"必须填写","")'
$orderWs.Range("C3").Formula  = '=IF(LEN(B3) =0,"必须填写","")'
$orderWs.Range("C4").Formula  = '=IF(LEN(B4) =0,"必须选择","")'
$orderWs.Range("C5").Formula  = '=IF(LEN(B5) =0,"必须填写","")'
$orderWs.Range("C6").Formula  = '=IF(LEN(B6) =0,"必须选择","")'
$orderWs.Range("C7").Formula  = '=IF(LEN(B7) =0,"必须填写","")'
$orderWs.Range("C8").Formula  = '=IF(LEN(B8) =0,"必须填写","")'
$orderWs.Range("C9").Formula  = '=IF(LEN(B9) =0,"必须选择","")'
$orderWs.Range("C11").Formula = '=IF(LEN(B11) =0,"必须选择","")'

# Move the active selection on "order" from C7 to B1.
$orderWs.Activate()
$orderWs.Range("B1").Select()

# --- Sheet "add": move the active selection from A1:XFD1048576 to A41 -----
$addWs = $wb.Worksheets.Item("add")
$addWs.Activate()
$addWs.Range("A41").Select()

# Restore "order" as the selected/active tab.
$orderWs.Activate()
